# Update "想去人数" (interest count) figures that were refreshed when the
# site's generated data was rebuilt at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 16490
$ws1.Range("F4").Value = 741
$ws1.Range("F5").Value = 254
$ws1.Range("F6").Value = 708
$ws1.Range("F7").Value = 1754
$ws1.Range("F8").Value = 165

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 16490
$ws4.Range("F4").Value = 741
$ws4.Range("F5").Value = 254
$ws4.Range("F8").Value = 709
$ws4.Range("F9").Value = 1754
$ws4.Range("F11").Value = 165
